# ---------------------------------------------------------------------------
# Applies the "properties -> xlsx" conversion edit:
#   * rename the sheet Sheet1 -> messages
#   * fix the German greeting string (drop "zurück")
#   * swap the user.manager.title / user.client.title rows (client now row 5,
#     manager now row 6)
#   * move the selection to A5
#   * touch row 7 (new blank formatted row under the table)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet.
$ws.Name = "messages"

# 2. Correct the German home-greetings translation.
$ws.Range("B2").Value = "Hallo und willkommen!"

# 3. Swap the "manager" and "client" rows (row 5 <-> row 6), keeping the
#    rest of the table untouched.
$managerRow = @($ws.Range("A5").Value2, $ws.Range("B5").Value2, $ws.Range("C5").Value2, $ws.Range("D5").Value2)
$clientRow  = @($ws.Range("A6").Value2, $ws.Range("B6").Value2, $ws.Range("C6").Value2, $ws.Range("D6").Value2)

$ws.Range("A5").Value = $clientRow[0]
$ws.Range("B5").Value = $clientRow[1]
$ws.Range("C5").Value = $clientRow[2]
$ws.Range("D5").Value = $clientRow[3]

$ws.Range("A6").Value = $managerRow[0]
$ws.Range("B6").Value = $managerRow[1]
$ws.Range("C6").Value = $managerRow[2]
$ws.Range("D6").Value = $managerRow[3]

# 4. New active selection cell.
$ws.Range("A5").Select() | Out-Null

# 5. The re-saved workbook recomputes the default row height for the header
#    rows (12.8 -> 12.75 pt); rows 5/6 keep the original 12.8pt height.
$ws.Rows.Item(1).RowHeight = 12.75
$ws.Rows.Item(2).RowHeight = 12.75
$ws.Rows.Item(3).RowHeight = 12.75
$ws.Rows.Item(4).RowHeight = 12.75

# 6. Give row 7 the same row height as the original table rows so it shows
#    up as a formatted-but-empty row below the data, like after the
#    conversion round-trip.
$ws.Rows.Item(7).RowHeight = 12.8
